# Weekly update: insert 3 new data rows for "Agrícola del Norte S.A. de Arica" - Cebolla
# at the top of this subgroup's block (row 577), pushing the existing rows (577-686) down
# by 3 (to 580-689). The new rows carry a more recent reporting date (44617).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 577 (existing row 577 onward shift down to 580 onward)
$ws.Range("A577:R579").EntireRow.Insert()

# --- Row 577 : 1a (cosecha) ---
$ws.Cells.Item(577, 1).Value  = 1
$ws.Cells.Item(577, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(577, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(577, 4).Value  = 44617
$ws.Cells.Item(577, 5).Value  = 15
$ws.Cells.Item(577, 6).Value  = 100112004
$ws.Cells.Item(577, 7).Value  = "Cebolla"
$ws.Cells.Item(577, 8).Value  = "Sin especificar"
$ws.Cells.Item(577, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(577, 10).Value = 1000
$ws.Cells.Item(577, 11).Value = 2500
$ws.Cells.Item(577, 12).Value = 3000
$ws.Cells.Item(577, 13).Value = 2750
$ws.Cells.Item(577, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(577, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(577, 16).Value = 153
$ws.Cells.Item(577, 17).Value = 18
$ws.Cells.Item(577, 18).Value = "Hortaliza"

# --- Row 578 : 2a (cosecha) ---
$ws.Cells.Item(578, 1).Value  = 1
$ws.Cells.Item(578, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(578, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(578, 4).Value  = 44617
$ws.Cells.Item(578, 5).Value  = 15
$ws.Cells.Item(578, 6).Value  = 100112004
$ws.Cells.Item(578, 7).Value  = "Cebolla"
$ws.Cells.Item(578, 8).Value  = "Sin especificar"
$ws.Cells.Item(578, 9).Value  = "2a (cosecha)"
$ws.Cells.Item(578, 10).Value = 1300
$ws.Cells.Item(578, 11).Value = 2000
$ws.Cells.Item(578, 12).Value = 2500
$ws.Cells.Item(578, 13).Value = 2250
$ws.Cells.Item(578, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(578, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(578, 16).Value = 125
$ws.Cells.Item(578, 17).Value = 18
$ws.Cells.Item(578, 18).Value = "Hortaliza"

# --- Row 579 : 3a (cosecha) ---
$ws.Cells.Item(579, 1).Value  = 1
$ws.Cells.Item(579, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(579, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(579, 4).Value  = 44617
$ws.Cells.Item(579, 5).Value  = 15
$ws.Cells.Item(579, 6).Value  = 100112004
$ws.Cells.Item(579, 7).Value  = "Cebolla"
$ws.Cells.Item(579, 8).Value  = "Sin especificar"
$ws.Cells.Item(579, 9).Value  = "3a (cosecha)"
$ws.Cells.Item(579, 10).Value = 1500
$ws.Cells.Item(579, 11).Value = 1500
$ws.Cells.Item(579, 12).Value = 2000
$ws.Cells.Item(579, 13).Value = 1750
$ws.Cells.Item(579, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(579, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(579, 16).Value = 97
$ws.Cells.Item(579, 17).Value = 18
$ws.Cells.Item(579, 18).Value = "Hortaliza"
